$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a cell to hold the given string as literal text (avoids Excel
# auto-converting number-like strings), then restore the default/Normal
# style so no stray formatting is left behind on the cell.
function Set-TextValue($ref, $value) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "64.727.10"
Set-TextValue "E2" "  +5.24%  "
Set-TextValue "D3" "3.104.52"
Set-TextValue "E3" "  +3.49%  "
Set-TextValue "E4" "  -0.05%  "
Set-TextValue "D5" "559.97"
Set-TextValue "E5" "  +2.50%  "
Set-TextValue "D6" "144.38"
Set-TextValue "E6" "  +10.41%  "
Set-TextValue "D7" "0.999"
Set-TextValue "D8" "3.099.04"
Set-TextValue "E8" "  +3.44%  "
Set-TextValue "E9" "  +2.31%  "
Set-TextValue "E10" "  +19.53%  "
Set-TextValue "E11" "  +5.14%  "
Set-TextValue "D12" "0.463"
Set-TextValue "E12" "  +4.30%  "
Set-TextValue "E13" "  +4.53%  "
Set-TextValue "D14" "35.45"
Set-TextValue "E14" "  +4.32%  "
Set-TextValue "D15" "3.607.05"
Set-TextValue "E15" "  +3.55%  "
Set-TextValue "D16" "64.697.62"
Set-TextValue "E16" "  +4.95%  "
Set-TextValue "D17" "3.108.14"
Set-TextValue "E17" "  +3.70%  "
Set-TextValue "E18" "  -0.43%  "
Set-TextValue "D19" "6.80"
Set-TextValue "E19" "  +3.05%  "
Set-TextValue "D20" "484.33"
Set-TextValue "E20" "  +0.89%  "
Set-TextValue "D21" "13.85"
Set-TextValue "E21" "  +5.31%  "
Set-TextValue "D22" "7.67"
Set-TextValue "E22" "  +10.23%  "
Set-TextValue "D23" "0.676"
Set-TextValue "E23" "  +2.02%  "
Set-TextValue "D24" "13.31"
Set-TextValue "E24" "  +10.92%  "
Set-TextValue "D25" "81.05"
Set-TextValue "E25" "  +0.50%  "
Set-TextValue "D26" "0.999"
Set-TextValue "E26" "  +0.00%  "
Set-TextValue "E27" "  +4.15%  "
Set-TextValue "D28" "8.11"
Set-TextValue "E28" "  +6.08%  "
Set-TextValue "D29" "2.08"
Set-TextValue "E29" "  +9.24%  "
Set-TextValue "E30" "  +0.02%  "
Set-TextValue "E31" "  +2.68%  "
Set-TextValue "E32" "  +3.78%  "
Set-TextValue "D33" "2.47"
Set-TextValue "E33" "  +5.98%  "
Set-TextValue "E34" "  +4.31%  "
Set-TextValue "B35" "Filecoin"
Set-TextValue "C35" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D35" "6.19"
Set-TextValue "E35" "  +5.97%  "
Set-TextValue "B36" "OKB"
Set-TextValue "C36" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D36" "55.31"
Set-TextValue "E36" "  +0.75%  "
Set-TextValue "D37" "466.68"
Set-TextValue "E37" "  +4.46%  "
Set-TextValue "D38" "0.0410"
Set-TextValue "E38" "  +7.70%  "
Set-TextValue "D39" "0.0830"
Set-TextValue "E39" "  +5.04%  "
Set-TextValue "D40" "3.019.53"
Set-TextValue "E40" "  -3.35%  "
Set-TextValue "E41" "  +1.46%  "
Set-TextValue "D42" "8.30"
Set-TextValue "E42" "  +2.92%  "
Set-TextValue "D43" "2.73"
Set-TextValue "E43" "  +15.94%  "
Set-TextValue "D44" "28.62"
Set-TextValue "E44" "  +12.20%  "
Set-TextValue "E45" "  +8.50%  "
Set-TextValue "E46" "  -0.01%  "
Set-TextValue "E47" "  +8.42%  "
Set-TextValue "D48" "0.112"
Set-TextValue "E48" "  +4.48%  "
Set-TextValue "D49" "118.83"
Set-TextValue "E49" "  +3.53%  "
Set-TextValue "E50" "  +6.97%  "
Set-TextValue "D51" "2.08"
Set-TextValue "E51" "  +3.20%  "
